# Fix column headers and add the new "PO Forecast" sheet with forecast data.

$wb = $excel.ActiveWorkbook

# --- Rename the "Requested quantity" headers on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match the bold/centered header style already used on the other sheets
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$forecastData = @(
    @(45025.99999999999,30,-12.89637330945477,69.07206842611673),
    @(45039.99999999999,30,-13.07572996994117,71.10833187444037),
    @(45046.99999999999,30,-9.518361272783023,74.8864761656195),
    @(45053.99999999999,31,-10.14030430621499,69.61993413179903),
    @(45060.99999999999,31,-7.098919186671057,72.03556277956301),
    @(45067.99999999999,31,-9.722744604708195,69.91790726913166),
    @(45081.99999999999,31,-10.12300040833935,71.64756400690683),
    @(45088.99999999999,31,-9.449143259210389,71.71203999034165),
    @(45116.99999999999,32,-5.996677324757979,72.02131628876222),
    @(45123.99999999999,32,-9.41341703005145,74.59961594782385),
    @(45137.99999999999,32,-10.23651858569512,72.13174726597984),
    @(45144.99999999999,32,-10.08505783821591,74.27403076897691),
    @(45151.99999999999,33,-10.4201566052212,71.29645586233154),
    @(45158.99999999999,33,-7.862537409194272,74.71676966545078),
    @(45165.99999999999,33,-9.52299264119155,73.58194374553234),
    @(45179.99999999999,33,-8.85918268545975,73.67124630429876),
    @(45186.99999999999,33,-7.064022125473108,72.6178307299239),
    @(45193.99999999999,34,-7.215497664935226,76.02562385360139),
    @(45200.99999999999,34,-6.430658325627443,74.5470385928415),
    @(45207.99999999999,34,-7.19988722539273,74.27635055792035),
    @(45214.99999999999,34,-10.21795678937115,76.8175671218674),
    @(45221.99999999999,34,-3.904591134033498,76.00286899071182),
    @(45228.99999999999,34,-7.84161407428102,73.98249236978884),
    @(45235.99999999999,34,-2.890514857209525,75.48752250572805),
    @(45242.99999999999,35,-6.499118741967825,78.19128312011526),
    @(45249.99999999999,35,-7.16640521983622,72.95234217212622),
    @(45256.99999999999,35,-4.532933440803401,75.81392753630419),
    @(45263.99999999999,35,-7.299870129033637,73.10130230623768),
    @(45270.99999999999,35,-4.143916161477298,76.45400805048185),
    @(45277.99999999999,35,-7.355239229619038,72.61460648375494)
)

for ($i = 0; $i -lt $forecastData.Length; $i++) {
    $row = $i + 2
    $rec = $forecastData[$i]

    $dsCell = $wsForecast.Cells.Item($row, 1)
    $dsCell.Value = $rec[0]
    $dsCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $wsForecast.Cells.Item($row, 2).Value = $rec[1]
    $wsForecast.Cells.Item($row, 3).Value = $rec[2]
    $wsForecast.Cells.Item($row, 4).Value = $rec[3]
}

$wsForecast.Range("A1").Select()
